$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 213, shifting existing rows 213-262 down to 214-263
$ws.Rows.Item(213).Insert()

# Populate the newly inserted row 213 with its data
$ws.Range("A213").Value = 5
$ws.Range("B213").Value = "Macroferia Regional de Talca"
$ws.Range("C213").Value = "Maule"
$ws.Range("D213").Value = 44511
$ws.Range("E213").Value = 7
$ws.Range("F213").Value = 100112043
$ws.Range("G213").Value = "Pepino ensalada"
$ws.Range("H213").Value = "Sin especificar"
$ws.Range("I213").Value = "Primera"
$ws.Range("J213").Value = 400
$ws.Range("K213").Value = 8000
$ws.Range("L213").Value = 8000
$ws.Range("M213").Value = 8000
$ws.Range("N213").Value = "$/caja 80 unidades"
$ws.Range("O213").Value = "Región del Maule"
$ws.Range("P213").Value = 100
$ws.Range("Q213").Value = 80
$ws.Range("R213").Value = "Hortaliza"
